$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1233387
$ws.Range("C2").Value = 1233387
$ws.Range("D2").Value = 1233387
$ws.Range("E2").Value = 1233387
$ws.Range("F2").Value = 880588
$ws.Range("G2").Value = 880588
$ws.Range("H2").Value = 1177787
$ws.Range("I2").Value = 1177787
$ws.Range("J2").Value = 1180865
$ws.Range("K2").Value = 1180870
$ws.Range("L2").Value = 1180870
$ws.Range("M2").Value = 1180870
$ws.Range("N2").Value = 1180870
$ws.Range("O2").Value = 790105
$ws.Range("P2").Value = 1186843
$ws.Range("Q2").Value = 1186843
$ws.Range("R2").Value = 1163096
$ws.Range("S2").Value = 1163096
$ws.Range("T2").Value = 1172439
$ws.Range("U2").Value = 737192
$ws.Range("Z2").Value = 84055
$ws.Range("AA2").Value = 84055
$ws.Range("AB2").Value = 99501
$ws.Range("AC2").Value = 99501
$ws.Range("AD2").Value = 83454
$ws.Range("AE2").Value = 83454
$ws.Range("C3").Value = 0.006434314614958647
$ws.Range("D3").Value = 5.187337226677433
$ws.Range("E3").Value = 0.4999808656974656
$ws.Range("F3").Value = 76.45911413737184
$ws.Range("G3").Value = 5.793132168505588
$ws.Range("H3").Value = 1.527495930928089
$ws.Range("I3").Value = 41.4027450464303
$ws.Range("J3").Value = 29.9279669699754
$ws.Range("K3").Value = 49.05968802662446
$ws.Range("L3").Value = 44.91276338631688
$ws.Range("M3").Value = 0.01776111680371252
$ws.Range("N3").Value = 44.38360276745112
$ws.Range("O3").Value = 30.99938362622689
$ws.Range("P3").Value = 6.872355265186717
$ws.Range("Q3").Value = -11.22799637357258
$ws.Range("R3").Value = 120.8834456485105
$ws.Range("S3").Value = 121.3917882100876
$ws.Range("T3").Value = 0.2172730913932409
$ws.Range("U3").Value = 0.2415331880975377
$ws.Range("Z3").Value = 51.11499518172625
$ws.Range("AB3").Value = 0.1856417523441976
$ws.Range("AC3").Value = 0.7426145465874715
$ws.Range("AD3").Value = -1.308756560500395
$ws.Range("AE3").Value = -1.308756560500395
$ws.Range("C4").Value = 0.1602993415109999
$ws.Range("D4").Value = 0.4663718094485888
$ws.Range("E4").Value = 2.329541842772553
$ws.Range("F4").Value = 18.85408367187142
$ws.Range("G4").Value = 1.379365628991676
$ws.Range("H4").Value = 9.617589749684138
$ws.Range("I4").Value = 295.4156733617752
$ws.Range("J4").Value = 1.133147246058553
$ws.Range("K4").Value = 16.67047149207881
$ws.Range("L4").Value = 39.64843055222235
$ws.Range("M4").Value = 0.08326084214932858
$ws.Range("N4").Value = 15.07719484408269
$ws.Range("O4").Value = 11.99670803555847
$ws.Range("P4").Value = 0.4593285532382717
$ws.Range("Q4").Value = 38.99206074345924
$ws.Range("R4").Value = 139.4460658705729
$ws.Range("S4").Value = 439.9277863798377
$ws.Range("T4").Value = 2.026561153747536
$ws.Range("U4").Value = 1.796087800601487
$ws.Range("Z4").Value = 49.24392449555241
$ws.Range("AB4").Value = 0.224440663687808
$ws.Range("AC4").Value = 0.897674536300694
$ws.Range("AD4").Value = 0.5557555302600938
$ws.Range("AE4").Value = 0.5557555302600938
$ws.Range("F6").Value = 59.6
$ws.Range("G6").Value = 4.58
$ws.Range("H6").Value = -1.88
$ws.Range("I6").Value = 30.34
$ws.Range("J6").Value = 29.181
$ws.Range("K6").Value = 45
$ws.Range("L6").Value = 41.5
$ws.Range("N6").Value = 41
$ws.Range("Q6").Value = -17.7
$ws.Range("R6").Value = 33.5
$ws.Range("S6").Value = 33.5
$ws.Range("T6").Value = -0.014
$ws.Range("U6").Value = -0.008999999999999999
$ws.Range("Z6").Value = 27.85
$ws.Range("AB6").Value = 0.06
$ws.Range("AC6").Value = 0.24
$ws.Range("AD6").Value = -1.32
$ws.Range("AE6").Value = -1.32
$ws.Range("F7").Value = 84.90000000000001
$ws.Range("G7").Value = 6.42
$ws.Range("H7").Value = 2.49
$ws.Range("I7").Value = 42.2
$ws.Range("J7").Value = 30.205
$ws.Range("K7").Value = 55.7
$ws.Range("L7").Value = 51.2
$ws.Range("P7").Value = 6.94
$ws.Range("Q7").Value = -14.1
$ws.Range("R7").Value = 50.1
$ws.Range("S7").Value = 50.1
$ws.Range("T7").Value = 0.012
$ws.Range("U7").Value = 0.02
$ws.Range("Z7").Value = 39.58
$ws.Range("AB7").Value = 0.14
$ws.Range("AC7").Value = 0.54
$ws.Range("AD7").Value = -1.24
$ws.Range("AE7").Value = -1.24
$ws.Range("F8").Value = 91
$ws.Range("G8").Value = 6.83
$ws.Range("H8").Value = 5.28
$ws.Range("I8").Value = 49.41
$ws.Range("J8").Value = 30.501
$ws.Range("Q8").Value = 11.7
$ws.Range("R8").Value = 235.2
$ws.Range("S8").Value = 235.2
$ws.Range("U8").Value = 0.047
$ws.Range("Z8").Value = 60.07
$ws.Range("AB8").Value = 0.3
$ws.Range("AC8").Value = 1.19
$ws.Range("AD8").Value = -1.12
$ws.Range("AE8").Value = -1.12
